# Project Sample Project is saved. Update Rules!B11 from "R40" to the
# text value "1" (kept as text, i.e. shared-string type, not a number).
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel's type-inference
# treat the numeric-looking text as a number (and typing it with a leading
# apostrophe marks the cell with a "number stored as text" quote-prefix,
# which changes its style). To land on genuine text without perturbing the
# cell's style, stage the text in a scratch cell, then Copy/PasteSpecial
# (values only) it into the target cell, and clean the scratch cell up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.Value = "'1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)

$scratch.Clear()
